# Atualizadas durante testes de atuma~
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows that had placeholder 0 values ---
$ws.Cells.Item(636, 3).Value = 806.228      # C636 - MAR 01/11/2023
$ws.Cells.Item(637, 3).Value = 798.898      # C637 - MAR 01/12/2023
$ws.Cells.Item(648, 3).Value = 36945.207    # C648 - TERRA 01/11/2023
$ws.Cells.Item(649, 3).Value = 39631.933    # C649 - TERRA 01/12/2023

# --- Append new MAR rows for 2024 (rows 650-661) ---
$marData = @(
    @("01/01/2024", 747.532),
    @("01/02/2024", 721.506),
    @("01/03/2024", 697.162),
    @("01/04/2024", 658.6369999999999),
    @("01/05/2024", 0),
    @("01/06/2024", 0),
    @("01/07/2024", 0),
    @("01/08/2024", 0),
    @("01/09/2024", 0),
    @("01/10/2024", 0),
    @("01/11/2024", 0),
    @("01/12/2024", 0)
)

$r = 650
foreach ($row in $marData) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.Style = "Normal"
    $ws.Cells.Item($r, 2).Value = "MAR"
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = "SERGIPE"
    $r = $r + 1
}

# --- Append new TERRA rows for 2024 (rows 662-673) ---
$terraData = @(
    @("01/01/2024", 42759.883),
    @("01/02/2024", 40222.657),
    @("01/03/2024", 44130.078),
    @("01/04/2024", 37515.475),
    @("01/05/2024", 0),
    @("01/06/2024", 0),
    @("01/07/2024", 0),
    @("01/08/2024", 0),
    @("01/09/2024", 0),
    @("01/10/2024", 0),
    @("01/11/2024", 0),
    @("01/12/2024", 0)
)

foreach ($row in $terraData) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $row[0]
    $cellA.Style = "Normal"
    $ws.Cells.Item($r, 2).Value = "TERRA"
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = "SERGIPE"
    $r = $r + 1
}
